$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.113.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.937.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9966"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.73"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9972"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4803"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2899"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06810"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.82"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "105.00"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07812"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.942.28"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.306"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6865"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "295.85"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +8.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.985.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007631"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.94"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9991"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.539"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.174.52"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.4665"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.36%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9967"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.417"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.583"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.39"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.86"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.123"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.395"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1016"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.648"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.532"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.361"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04841"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7412"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.132"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.723"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01964"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.550"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.635"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "77.08"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.047"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8752"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4373"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9981"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.025.13"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.51%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.595"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1216"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.109"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.45%  "
